# Added check_in/add_new methods and started check_out
# - Condensed Inventory: two new wines added (wine_id 5 "Syrah" / wine_id 6 "Zinfandel")
# - Expanded Inventory: matching per-bottle rows with location/date_in tracking

$wb = $excel.ActiveWorkbook
$condensed = $wb.Worksheets.Item("Condensed Inventory")
$expanded  = $wb.Worksheets.Item("Expanded Inventory")

function Set-Text($ws, $row, $col, $text) {
    # Force text so values like "$35" aren't auto-parsed as currency numbers,
    # then clear the formatting byproduct so the cell keeps the default style.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Condensed Inventory: row 3 (wine_id 5, Syrah) ---
$condensed.Cells.Item(3, 1).Value = 5
Set-Text $condensed 3 3 "Burnt Bridge Cellars"
Set-Text $condensed 3 4 "Walla Walla"
Set-Text $condensed 3 6 "Syrah"
$condensed.Cells.Item(3, 7).Value = 2014
Set-Text $condensed 3 8 "Table"
Set-Text $condensed 3 9 '$35'
Set-Text $condensed 3 10 '$35'
Set-Text $condensed 3 12 "Standard (750 mL)"
$condensed.Cells.Item(3, 13).Value = 1

# --- Condensed Inventory: row 4 (wine_id 6, Zinfandel) ---
$condensed.Cells.Item(4, 1).Value = 6
Set-Text $condensed 4 3 "Turly"
Set-Text $condensed 4 4 "Walla Walla"
Set-Text $condensed 4 6 "Zinfandel"
$condensed.Cells.Item(4, 7).Value = 2011
Set-Text $condensed 4 8 "Table"
Set-Text $condensed 4 9 '$40'
Set-Text $condensed 4 10 '$25'
Set-Text $condensed 4 11 "Young vines"
Set-Text $condensed 4 12 "Standard (750 mL)"
$condensed.Cells.Item(4, 13).Value = 2

# --- Expanded Inventory: row 4 (wine_id 5, Syrah, add_new @ C6) ---
$expanded.Cells.Item(4, 1).Value = 5
Set-Text $expanded 4 3 "Burnt Bridge Cellars"
Set-Text $expanded 4 4 "Walla Walla"
Set-Text $expanded 4 6 "Syrah"
$expanded.Cells.Item(4, 7).Value = 2014
Set-Text $expanded 4 8 "Table"
Set-Text $expanded 4 9 '$35'
Set-Text $expanded 4 10 '$35'
Set-Text $expanded 4 12 "Standard (750 mL)"
Set-Text $expanded 4 13 "C6"
Set-Text $expanded 4 14 "2019-01-22 16:35:03"

# --- Expanded Inventory: row 5 (wine_id 6, Zinfandel, check_in #1 @ C6) ---
$expanded.Cells.Item(5, 1).Value = 6
Set-Text $expanded 5 3 "Turly"
Set-Text $expanded 5 4 "Walla Walla"
Set-Text $expanded 5 6 "Zinfandel"
$expanded.Cells.Item(5, 7).Value = 2011
Set-Text $expanded 5 8 "Table"
Set-Text $expanded 5 9 '$40'
Set-Text $expanded 5 10 '$25'
Set-Text $expanded 5 11 "Young vines"
Set-Text $expanded 5 12 "Standard (750 mL)"
Set-Text $expanded 5 13 "C6"
Set-Text $expanded 5 14 "2019-01-22 17:49:41"

# --- Expanded Inventory: row 6 (wine_id 6, Zinfandel, check_in #2 @ C6) ---
$expanded.Cells.Item(6, 1).Value = 6
Set-Text $expanded 6 3 "Turly"
Set-Text $expanded 6 4 "Walla Walla"
Set-Text $expanded 6 6 "Zinfandel"
$expanded.Cells.Item(6, 7).Value = 2011
Set-Text $expanded 6 8 "Table"
Set-Text $expanded 6 9 '$40'
Set-Text $expanded 6 10 '$25'
Set-Text $expanded 6 11 "Young vines"
Set-Text $expanded 6 12 "Standard (750 mL)"
Set-Text $expanded 6 13 "C6"
Set-Text $expanded 6 14 "2019-01-22 17:53:58"
